$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 and 14: Coin name and Link swap
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"

# Price (Column D) updates -- prefix with apostrophe to force text and avoid numeric coercion
$ws.Range("D2").Value = "'25.780.02"
$ws.Range("D3").Value = "'1.626.99"
$ws.Range("D4").Value = "'1.002"
$ws.Range("D5").Value = "'215.08"
$ws.Range("D6").Value = "'0.5066"
$ws.Range("D7").Value = "'1.003"
$ws.Range("D8").Value = "'0.2577"
$ws.Range("D9").Value = "'0.06445"
$ws.Range("D10").Value = "'19.40"
$ws.Range("D11").Value = "'0.07787"
$ws.Range("D12").Value = "'4.255"
$ws.Range("D13").Value = "'1.853.03"
$ws.Range("D14").Value = "'1.622.29"
$ws.Range("D15").Value = "'0.5557"
$ws.Range("D16").Value = "'63.03"
$ws.Range("D17").Value = "'0.0₅7535"
$ws.Range("D18").Value = "'25.795.95"
$ws.Range("D19").Value = "'1.003"
$ws.Range("D20").Value = "'193.63"
$ws.Range("D21").Value = "'4.293"
$ws.Range("D22").Value = "'9.805"
$ws.Range("D23").Value = "'6.001"
$ws.Range("D24").Value = "'1.003"
$ws.Range("D25").Value = "'1.809"
$ws.Range("D26").Value = "'140.44"
$ws.Range("D27").Value = "'0.1259"
$ws.Range("D28").Value = "'6.714"
$ws.Range("D29").Value = "'15.38"
$ws.Range("D30").Value = "'1.235"
$ws.Range("D31").Value = "'0.04869"
$ws.Range("D32").Value = "'3.269"
$ws.Range("D33").Value = "'3.172"
$ws.Range("D34").Value = "'1.553"
$ws.Range("D35").Value = "'2.372"
$ws.Range("D36").Value = "'0.8928"
$ws.Range("D37").Value = "'2.569"
$ws.Range("D38").Value = "'1.132.59"
$ws.Range("D39").Value = "'0.5457"
$ws.Range("D40").Value = "'0.01553"
$ws.Range("D41").Value = "'0.9919"
$ws.Range("D42").Value = "'5.561"
$ws.Range("D43").Value = "'0.7939"
$ws.Range("D44").Value = "'97.24"
$ws.Range("D45").Value = "'1.781.49"
$ws.Range("D46").Value = "'0.0₈112"
$ws.Range("D47").Value = "'0.4432"
$ws.Range("D48").Value = "'55.00"
$ws.Range("D49").Value = "'0.05049"
$ws.Range("D50").Value = "'7.591"

# Volume(1h) (Column E) updates
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +1.41%  "
$ws.Range("E10").Value = "  -2.26%  "
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("E16").Value = "  -2.04%  "
$ws.Range("E17").Value = "  -2.68%  "
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("E20").Value = "  -1.02%  "
$ws.Range("E21").Value = "  -3.17%  "
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("E23").Value = "  -1.85%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  -4.41%  "
$ws.Range("E26").Value = "  -2.13%  "
$ws.Range("E27").Value = "  +1.49%  "
$ws.Range("E28").Value = "  -2.49%  "
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("E33").Value = "  -1.34%  "
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("E36").Value = "  -2.57%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  +3.79%  "
$ws.Range("E39").Value = "  -1.58%  "
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("E43").Value = "  -1.52%  "
$ws.Range("E44").Value = "  -1.69%  "
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("E46").Value = "  -6.58%  "
$ws.Range("E47").Value = "  -2.30%  "
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("E49").Value = "  -3.24%  "
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("E51").Value = "  +0.26%  "
